# 20200303 Update for examnple excel
#
# - CMD!B6 had a stray trailing space in the "read" command text; fix it to "read".
# - The CMD sheet becomes the active/selected sheet (instead of PORT), with the
#   current selection at B6 (instead of C13), and the PORT sheet's prior scroll
#   position (topLeftCell=C1) / tab-selected flag are cleared since PORT is no
#   longer the front-most sheet.

$wb = $excel.ActiveWorkbook

$cmdSheet = $wb.Worksheets.Item("CMD")

# Fix the typo'd trailing space in the SPI "read " command.
$cmdSheet.Range("B6").Value = "read"

# Make CMD the active sheet/tab, with B6 selected (matches the saved view state).
$cmdSheet.Activate()
$cmdSheet.Range("B6").Select()
